$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 3867.6667
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3867.6667
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 11603.0001
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -11841.0001
$ws.Range("H60").Value = 3867.6667
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 3867.6667
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 11603.0001
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -12571.0001
$ws.Range("H137").Value = 953.6087
$ws.Range("I137").Value = 1333.3334
$ws.Range("J137").Value = 896.65
$ws.Range("K137").Value = 4000.0002
$ws.Range("L137").Value = 2689.95
$ws.Range("M137").Value = -1450.0002
$ws.Range("N137").Value = -7789.95
$ws.Range("H138").Value = 2031.1616
$ws.Range("I138").Value = 921.5641000000001
$ws.Range("J138").Value = 2752.4
$ws.Range("K138").Value = 2764.6923
$ws.Range("L138").Value = 8257.200000000001
$ws.Range("M138").Value = 2375.3077
$ws.Range("N138").Value = -18537.2

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2483.158
$ws.Range("I88").Value = 2044.6154
$ws.Range("J88").Value = 3433.3333
$ws.Range("K88").Value = 2044.6154
$ws.Range("L88").Value = 3433.3333
$ws.Range("M88").Value = -1638.6154
$ws.Range("N88").Value = -4245.3333
$ws.Range("H91").Value = 2483.158
$ws.Range("I91").Value = 2044.6154
$ws.Range("J91").Value = 3433.3333
$ws.Range("K91").Value = 2044.6154
$ws.Range("L91").Value = 3433.3333
$ws.Range("M91").Value = -640.6153999999999
$ws.Range("N91").Value = -6241.3333
$ws.Range("H132").Value = 1127
$ws.Range("I132").Value = 999.23334
$ws.Range("J132").Value = 1552.8889
$ws.Range("K132").Value = 2997.70002
$ws.Range("L132").Value = 4658.6667
$ws.Range("M132").Value = -467.7000200000002
$ws.Range("N132").Value = -9718.6667

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2270.5151
$ws.Range("I86").Value = 2119.6553
$ws.Range("K86").Value = 2119.6553
$ws.Range("M86").Value = -996.6552999999999
$ws.Range("H89").Value = 2270.5151
$ws.Range("I89").Value = 2119.6553
$ws.Range("K89").Value = 10598.2765
$ws.Range("M89").Value = -4982.2765

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2145420.5
$ws.Range("I31").Value = 2502360
$ws.Range("K31").Value = 2502360
$ws.Range("M31").Value = -2502065
$ws.Range("H34").Value = 2145420.5
$ws.Range("I34").Value = 2502360
$ws.Range("K34").Value = 2502360
$ws.Range("M34").Value = -2502158
$ws.Range("H132").Value = 2210.3262
$ws.Range("I132").Value = 1436.4286
$ws.Range("K132").Value = 4309.2858
$ws.Range("M132").Value = -1779.2858
$ws.Range("H134").Value = 1409.28
$ws.Range("I134").Value = 1420.3125
$ws.Range("J134").Value = 1389.6666
$ws.Range("K134").Value = 4260.9375
$ws.Range("L134").Value = 4168.9998
$ws.Range("M134").Value = -1725.9375
$ws.Range("N134").Value = -9238.9998

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 797.871
$ws.Range("I5").Value = 900.1818
$ws.Range("J5").Value = 547.7778
$ws.Range("K5").Value = 2700.5454
$ws.Range("L5").Value = 1643.3334
$ws.Range("M5").Value = -2588.5454
$ws.Range("N5").Value = -1867.3334
$ws.Range("H107").Value = 519083.4
$ws.Range("I107").Value = 797
$ws.Range("J107").Value = 1111410.8
$ws.Range("K107").Value = 2391
$ws.Range("L107").Value = 3334232.4
$ws.Range("M107").Value = -471
$ws.Range("N107").Value = -3338072.4
$ws.Range("H126").Value = 63066.47
$ws.Range("I126").Value = 145118.58
$ws.Range("J126").Value = 5630
$ws.Range("K126").Value = 435355.74
$ws.Range("L126").Value = 16890
$ws.Range("M126").Value = -430415.74
$ws.Range("N126").Value = -26770
$ws.Range("H129").Value = 49897.855
$ws.Range("J129").Value = 86413.75
$ws.Range("L129").Value = 259241.25
$ws.Range("N129").Value = -269241.25
$ws.Range("H131").Value = 33805.38
$ws.Range("J131").Value = 20188.02
$ws.Range("L131").Value = 60564.06
$ws.Range("N131").Value = -70644.06
$ws.Range("H132").Value = 679.6
$ws.Range("I132").Value = 607.63635
$ws.Range("J132").Value = 877.5
$ws.Range("K132").Value = 5468.72715
$ws.Range("L132").Value = 7897.5
$ws.Range("M132").Value = -2938.72715
$ws.Range("N132").Value = -12957.5
$ws.Range("H133").Value = 4926
$ws.Range("I133").Value = 1576.6666
$ws.Range("J133").Value = 9950
$ws.Range("K133").Value = 4729.9998
$ws.Range("L133").Value = 29850
$ws.Range("M133").Value = 330.0002000000004
$ws.Range("N133").Value = -39970
$ws.Range("H134").Value = 2842.647
$ws.Range("I134").Value = 1642.5927
$ws.Range("J134").Value = 7471.4287
$ws.Range("K134").Value = 4927.7781
$ws.Range("L134").Value = 22414.2861
$ws.Range("M134").Value = 142.2219000000005
$ws.Range("N134").Value = -32554.2861
$ws.Range("H135").Value = 797.871
$ws.Range("I135").Value = 900.1818
$ws.Range("J135").Value = 547.7778
$ws.Range("K135").Value = 8101.6362
$ws.Range("L135").Value = 4930.000199999999
$ws.Range("M135").Value = -5566.6362
$ws.Range("N135").Value = -10000.0002
$ws.Range("H136").Value = 53159.5
$ws.Range("I136").Value = 101855
$ws.Range("J136").Value = 4464
$ws.Range("K136").Value = 305565
$ws.Range("L136").Value = 13392
$ws.Range("M136").Value = -300465
$ws.Range("N136").Value = -23592
$ws.Range("H137").Value = 45585.125
$ws.Range("I137").Value = 2247.647
$ws.Range("J137").Value = 150833.28
$ws.Range("K137").Value = 6742.941
$ws.Range("L137").Value = 452499.84
$ws.Range("M137").Value = -1642.941
$ws.Range("N137").Value = -462699.84
$ws.Range("H138").Value = 1545
$ws.Range("I138").Value = 1051.4286
$ws.Range("J138").Value = 5000
$ws.Range("K138").Value = 3154.2858
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = 1985.7142
$ws.Range("N138").Value = -25280
$ws.Range("H139").Value = 85148.586
$ws.Range("I139").Value = 92586.37
$ws.Range("J139").Value = 3333
$ws.Range("K139").Value = 277759.11
$ws.Range("L139").Value = 9999
$ws.Range("M139").Value = -272619.11
$ws.Range("N139").Value = -20279

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 19900
$ws.Range("J15").Value = 19900
$ws.Range("L15").Value = 19900
$ws.Range("N15").Value = -20476
$ws.Range("H70").Value = 4856.6
$ws.Range("I70").Value = 4844.3335
$ws.Range("J70").Value = 4875
$ws.Range("K70").Value = 4844.3335
$ws.Range("L70").Value = 4875
$ws.Range("M70").Value = -4574.3335
$ws.Range("N70").Value = -5415
$ws.Range("H73").Value = 4856.6
$ws.Range("I73").Value = 4844.3335
$ws.Range("J73").Value = 4875
$ws.Range("K73").Value = 4844.3335
$ws.Range("L73").Value = 4875
$ws.Range("M73").Value = -3908.3335
$ws.Range("N73").Value = -6747
$ws.Range("H81").Value = 19900
$ws.Range("J81").Value = 19900
$ws.Range("L81").Value = 19900
$ws.Range("N81").Value = -21896
$ws.Range("H84").Value = 19900
$ws.Range("J84").Value = 19900
$ws.Range("L84").Value = 59700
$ws.Range("N84").Value = -69684

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3088.2856
$ws.Range("I132").Value = 2998.5789
$ws.Range("J132").Value = 3277.6667
$ws.Range("K132").Value = 8995.736699999999
$ws.Range("L132").Value = 9833.000100000001
$ws.Range("M132").Value = -6465.736699999999
$ws.Range("N132").Value = -14893.0001
$ws.Range("H136").Value = 1671.2903
$ws.Range("I136").Value = 934.3396
$ws.Range("J136").Value = 6011.1113
$ws.Range("K136").Value = 2803.0188
$ws.Range("L136").Value = 18033.3339
$ws.Range("M136").Value = -253.0187999999998
$ws.Range("N136").Value = -23133.3339

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 10133.333
$ws.Range("I26").Value = 11900
$ws.Range("J26").Value = 9250
$ws.Range("K26").Value = 11900
$ws.Range("L26").Value = 9250
$ws.Range("M26").Value = -11607
$ws.Range("N26").Value = -9836
$ws.Range("H132").Value = 891.75
$ws.Range("I132").Value = 705.5454999999999
$ws.Range("K132").Value = 2116.6365
$ws.Range("M132").Value = 413.3635000000004
$ws.Range("H136").Value = 1043.0968
$ws.Range("I136").Value = 1093.48
$ws.Range("J136").Value = 833.1667
$ws.Range("K136").Value = 3280.44
$ws.Range("L136").Value = 2499.5001
$ws.Range("M136").Value = -730.4400000000001
$ws.Range("N136").Value = -7599.5001

Write-Host "Applied Bahamut_Profits market data refresh"